# 4.4 Caso de Teste - UC-41 Cadastrar produto
# Commit: "Correcoes nos diagramas e prototipos"
#
# Changes applied:
#  1) "BOTAO SELECIONAR MAQUINA"          -> "BOTAO SELECIONAR FORNECEDOR" (x2)
#  2) "...CAMPO MAQUINA DA TELA..."       -> "...CAMPO FORNECEDOR DA TELA..."
#  3) "TODAS OS FORNECEDORES..."          -> "TODOS OS FORNECEDORES..." (typo fix)
#     + the _GoBack bookmark is relocated from the trailing empty paragraph
#       (which is removed) into the middle of that same cell's text.
#  4) "SALVAR" / "BOTAO SALVAR" button captions -> "CADASTRAR" / "BOTAO CADASTRAR" (x5)
#  5) The results table's last two columns are widened: 1340 -> 1482, 1984 -> 1842

$d = $word.ActiveDocument

# --- 1) "BOTÃO SELECIONAR MÁQUINA" -> "BOTÃO SELECIONAR FORNECEDOR" (both occurrences) ---
$d.Content.Find.Execute("BOTÃO SELECIONAR MÁQUINA", $false, $true, $false, $false, $false, `
                         $true, 1, $false, "BOTÃO SELECIONAR FORNECEDOR", 2) | Out-Null

# --- 2) "...PREENCHENDO O CAMPO MÁQUINA DA TELA SOBREPOSTA." -> "...CAMPO FORNECEDOR..." ---
$d.Content.Find.Execute("CAMPO MÁQUINA DA TELA SOBREPOSTA", $false, $true, $false, $false, $false, `
                         $true, 1, $false, "CAMPO FORNECEDOR DA TELA SOBREPOSTA", 2) | Out-Null

# --- 3) Fix "TODAS OS FORNECEDORES..." -> "TODOS OS FORNECEDORES..." ---
$d.Content.Find.Execute("TODAS OS FORNECEDORES CADASTRADOS NO SISTEMA SÃO LISTADAS.", `
                         $false, $true, $false, $false, $false, `
                         $true, 1, $false, "TODOS OS FORNECEDORES CADASTRADOS NO SISTEMA SÃO LISTADAS.", 2) | Out-Null

# --- 4) "SALVAR" -> "CADASTRAR" on every button caption (5 occurrences, incl. "BOTÃO SALVAR") ---
for ($i = 0; $i -lt 5; $i++) {
    $d.Content.Find.Execute("SALVAR", $false, $true, $false, $false, $false, `
                             $true, 0, $false, "CADASTRAR", 1) | Out-Null
}

# --- 5) Move the _GoBack bookmark ---
# Remove it (and the otherwise-empty paragraph that hosts it) from its old spot...
if ($d.Bookmarks.Exists("_GoBack")) {
    $old = $d.Bookmarks("_GoBack")
    $oldPos = $old.Start
    $para = $d.Range($oldPos, $oldPos).Paragraphs(1)
    $pStart = $para.Range.Start
    $pEnd = $para.Range.End
    $d.Range($pStart, $pEnd).Delete() | Out-Null
}

# ...and re-create it inside "TODOS OS FORNECEDORES..." cell, splitting "TODOS" into
# "TODO" + bookmark + "S OS FORNECEDORES..." just like the tracked edit did.
$target = $d.Content
$target.Find.Execute("TODOS OS FORNECEDORES CADASTRADOS NO SISTEMA SÃO LISTADAS.", `
                      $false, $true, $false, $false, $false, `
                      $true, 0, $false, "", 0) | Out-Null
$newPos = $target.Start + 4
$d.Bookmarks.Add("_GoBack", $d.Range($newPos, $newPos)) | Out-Null

# --- 6) Widen the last two columns of the results table (8-column table): 1340->1482, 1984->1842 ---
for ($ti = 1; $ti -le $d.Tables.Count; $ti++) {
    $tbl = $d.Tables($ti)
    if ($tbl.Columns.Count -eq 8) {
        $tbl.Columns(7).Width = 74.1   # 1482 twips (points = twips / 20)
        $tbl.Columns(8).Width = 92.1   # 1842 twips
    }
}
